$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.470362067222595
$ws.Range("B1").Value = 1.725765109062195
$ws.Range("C1").Value = 1.641141414642334
$ws.Range("D1").Value = 1.53488028049469
$ws.Range("E1").Value = 1.082818269729614
